$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = "Aufgabe 5 Rekurrenz"
$ws.Range("B23").Value = "Random Neighbor selection"

$ws.Range("B23").Select()
